$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete row 2 (the duplicate header row). This shifts the old rows 3-10
# up to become rows 2-9 (each keeping its own values/types), and leaves a
# brand-new blank row 10 behind.
$ws.Rows.Item(2).Delete()

# Row 1 loses its bold/centered/bordered style, and A1's "Unnamed: 0" text
# is cleared. Assigning a bare quote keeps the cell typed as Text/empty
# (matching how the sheet stores blank string cells) instead of Excel's
# default blank-Number typing; ClearFormats then drops the style back to
# the sheet default.
$ws.Range("A1").Value = "'"
$ws.Range("A1:P1").ClearFormats()

# Row 10 is the freshly shifted-in blank row -- every cell defaulted to an
# untyped/Number blank, so restamp the whole row as empty Text cells (again
# using the quote-prefix trick) and clear the resulting quote-prefix style
# back to the sheet default.
$ws.Range("A10:P10").Value = "'"
$ws.Range("A10:P10").ClearFormats()

# Add a brand new, entirely empty row 11 underneath.
$ws.Range("A11:P11").Value = "'"
$ws.Range("A11:P11").ClearFormats()
